$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nifty")
$ws.Range("G1").Value = 43448
$ws.Range("H1").Value = 43451
$ws.Range("G2").Value = $ws.Range("H2").Value
$ws.Range("H2").Value = 10900.35
